$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, pushing old rows 65-94 down to 66-95
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with its data
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Vega Modelo de Temuco"
$ws.Range("C65").Value = "La Araucanía"
$ws.Range("D65").Value = 45141
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = 100112042
$ws.Range("G65").Value = "Locoto"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 80
$ws.Range("K65").Value = 2700
$ws.Range("L65").Value = 2700
$ws.Range("M65").Value = 2700
$ws.Range("N65").Value = "`$/kilo"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 2700
$ws.Range("Q65").Value = 1
$ws.Range("R65").Value = "Hortaliza"
